$d = $word.ActiveDocument

$replacements = @(
    "69÷8=8, 5", "27÷2=13, 1",
    "87÷3=29, 0", "28÷2=14, 0",
    "77÷2=38, 1", "82÷3=27, 1",
    "21÷2=10, 1", "88÷2=44, 0",
    "26÷2=13, 0", "31÷4=7, 3",
    "61÷8=7, 5", "36÷9=4, 0",
    "91÷3=30, 1", "31÷7=4, 3",
    "74÷6=12, 2", "29÷2=14, 1",
    "88÷5=17, 3", "72÷7=10, 2",
    "82÷5=16, 2", "83÷4=20, 3",
    "87÷9=9, 6", "72÷2=36, 0",
    "49÷9=5, 4", "63÷8=7, 7",
    "99÷9=11, 0", "41÷6=6, 5",
    "61÷7=8, 5", "34÷9=3, 7",
    "11÷8=1, 3", "46÷7=6, 4",
    "89÷5=17, 4", "33÷8=4, 1",
    "98÷2=49, 0", "61÷3=20, 1",
    "21÷2=10, 1", "41÷8=5, 1",
    "21÷4=5, 1", "40÷7=5, 5",
    "75÷9=8, 3", "26÷4=6, 2",
    "28÷2=14, 0", "63÷4=15, 3",
    "39÷7=5, 4", "96÷9=10, 6",
    "34÷5=6, 4", "21÷5=4, 1",
    "10÷2=5, 0", "66÷3=22, 0",
    "49÷3=16, 1", "75÷5=15, 0"
)

$table = $d.Tables.Item(1)
$idx = 0
for ($r = 1; $r -le $table.Rows.Count; $r++) {
    for ($c = 1; $c -le $table.Columns.Count; $c++) {
        $cell = $null
        try {
            $cell = $table.Cell($r, $c)
        } catch {
            $cell = $null
        }
        if ($cell -ne $null) {
            $rng = $cell.Range
            $cellText = $rng.Text
            if ($cellText -ne $null -and $cellText.Length -gt 2) {
                $old = $replacements[$idx * 2]
                $new = $replacements[$idx * 2 + 1]
                # Cell.Range spans the cell contents plus the trailing
                # end-of-cell mark (2 chars); target just the content.
                $contentRng = $d.Range($rng.Start, $rng.End - 1)
                if ($contentRng.Text -eq $old) {
                    $contentRng.Text = $new
                }
                $idx = $idx + 1
            }
        }
    }
}
